$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "29.555.29"
$ws.Range("D2").ClearFormats()
$ws.Range("E2").Value = "  +3.07%  "

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "1.603.96"
$ws.Range("D3").ClearFormats()
$ws.Range("E3").Value = "  +2.38%  "

$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "0.999"
$ws.Range("D4").ClearFormats()
$ws.Range("E4").Value = "  -0.04%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "212.64"
$ws.Range("D5").ClearFormats()
$ws.Range("E5").Value = "  +1.21%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "0.523"
$ws.Range("D6").ClearFormats()
$ws.Range("E6").Value = "  +7.00%  "

$ws.Range("E7").Value = "  -0.03%  "

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "26.95"
$ws.Range("D8").ClearFormats()
$ws.Range("E8").Value = "  +7.52%  "

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "43.45"
$ws.Range("D9").ClearFormats()

$ws.Range("E10").Value = "  +2.63%  "

$ws.Range("E11").Value = "  +2.76%  "

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.0910"
$ws.Range("D12").ClearFormats()
$ws.Range("E12").Value = "  +1.57%  "

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "1.833.14"
$ws.Range("D13").ClearFormats()
$ws.Range("E13").Value = "  +2.32%  "

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "1.630.10"
$ws.Range("D14").ClearFormats()
$ws.Range("E14").Value = "  +4.14%  "

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "29.568.51"
$ws.Range("D15").ClearFormats()
$ws.Range("E15").Value = "  +3.08%  "

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "0.537"
$ws.Range("D16").ClearFormats()
$ws.Range("E16").Value = "  +3.63%  "

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "3.72"
$ws.Range("D17").ClearFormats()
$ws.Range("E17").Value = "  +2.00%  "

$ws.Range("E18").Value = "  +3.40%  "

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "240.78"
$ws.Range("D19").ClearFormats()
$ws.Range("E19").Value = "  +4.94%  "

$ws.Range("E20").Value = "  +3.38%  "

$ws.Range("E21").Value = "  +1.82%  "

$ws.Range("E22").Value = "  +0.01%  "

$ws.Range("E23").Value = "  +1.53%  "

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "9.22"
$ws.Range("D24").ClearFormats()
$ws.Range("E24").Value = "  +2.11%  "

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "2.09"
$ws.Range("D25").ClearFormats()
$ws.Range("E25").Value = "  +0.70%  "

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "154.50"
$ws.Range("D26").ClearFormats()
$ws.Range("E26").Value = "  +1.69%  "

$ws.Range("E27").Value = "  +5.08%  "

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "15.26"
$ws.Range("D28").ClearFormats()
$ws.Range("E28").Value = "  +3.20%  "

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "6.40"
$ws.Range("D29").ClearFormats()
$ws.Range("E29").Value = "  +2.55%  "

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "0.998"
$ws.Range("D30").ClearFormats()
$ws.Range("E30").Value = "  -0.09%  "

$ws.Range("E31").Value = "  +2.71%  "

$ws.Range("E32").Value = "  +0.67%  "

$ws.Range("E33").Value = "  +1.64%  "

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "3.10"
$ws.Range("D34").ClearFormats()
$ws.Range("E34").Value = "  +3.66%  "

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "1.409.87"
$ws.Range("D35").ClearFormats()
$ws.Range("E35").Value = "  +1.12%  "

$ws.Range("E36").Value = "  +0.76%  "

$ws.Range("E37").Value = "  +4.14%  "

$ws.Range("E38").Value = "  +4.18%  "

$ws.Range("E39").Value = "  +0.09%  "

$ws.Range("E40").Value = "  +2.30%  "

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.539"
$ws.Range("D41").ClearFormats()
$ws.Range("E41").Value = "  +3.44%  "

$ws.Range("E42").Value = "  +0.88%  "

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.0490"
$ws.Range("D43").ClearFormats()
$ws.Range("E43").Value = "  +6.65%  "

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "53.15"
$ws.Range("D44").ClearFormats()
$ws.Range("E44").Value = "  +23.62%  "

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.798"
$ws.Range("D45").ClearFormats()
$ws.Range("E45").Value = "  +3.34%  "

$ws.Range("E46").Value = "  -0.04%  "

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "65.82"
$ws.Range("D47").ClearFormats()
$ws.Range("E47").Value = "  +2.89%  "

$ws.Range("E48").Value = "  +0.72%  "

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "1.744.87"
$ws.Range("D49").ClearFormats()
$ws.Range("E49").Value = "  +2.57%  "

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.857"
$ws.Range("D50").ClearFormats()
$ws.Range("E50").Value = "  -1.49%  "

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "86.56"
$ws.Range("D51").ClearFormats()
$ws.Range("E51").Value = "  +1.74%  "
